# Update stock-report rows whose item-code/qty/rate columns (B, D, E, F, G)
# were rotated/swapped between duplicate product rows in the source report.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B136").Value = 63902
$ws.Range("D136").Value = 32.02
$ws.Range("E136").Value = 34.04
$ws.Range("F136").Value = 2
$ws.Range("G136").Value = 64.04000000000001
$ws.Range("B137").Value = 48654
$ws.Range("D137").Value = 32.02
$ws.Range("E137").Value = 38.26
$ws.Range("F137").Value = -1
$ws.Range("G137").Value = -32.02
$ws.Range("B146").Value = 53925
$ws.Range("D146").Value = 66.44
$ws.Range("E146").Value = 79.37
$ws.Range("F146").Value = 1
$ws.Range("G146").Value = 66.44
$ws.Range("B147").Value = 64350
$ws.Range("D147").Value = 66.44
$ws.Range("E147").Value = 70.63
$ws.Range("F147").Value = 2
$ws.Range("G147").Value = 132.88
$ws.Range("B148").Value = 57756
$ws.Range("D148").Value = 66.44
$ws.Range("E148").Value = 79.37
$ws.Range("F148").Value = -100
$ws.Range("G148").Value = -6644
$ws.Range("B233").Value = 48719
$ws.Range("D233").Value = 295.75
$ws.Range("E233").Value = 353.35
$ws.Range("F233").Value = -81
$ws.Range("G233").Value = -23955.75
$ws.Range("B234").Value = 64979
$ws.Range("D234").Value = 295.75
$ws.Range("E234").Value = 314.41
$ws.Range("F234").Value = 25
$ws.Range("G234").Value = 7393.75
$ws.Range("B246").Value = 48706
$ws.Range("D246").Value = 33.3
$ws.Range("E246").Value = 39.8
$ws.Range("F246").Value = -144
$ws.Range("G246").Value = -4795.2
$ws.Range("B247").Value = 64973
$ws.Range("D247").Value = 33.3
$ws.Range("E247").Value = 35.4
$ws.Range("F247").Value = 99
$ws.Range("G247").Value = 3296.7
$ws.Range("B294").Value = 63571
$ws.Range("D294").Value = 143.48
$ws.Range("E294").Value = 152.53
$ws.Range("F294").Value = 8
$ws.Range("G294").Value = 1147.84
$ws.Range("B295").Value = 63531
$ws.Range("D295").Value = 143.48
$ws.Range("E295").Value = 152.53
$ws.Range("F295").Value = 80
$ws.Range("G295").Value = 11478.4
$ws.Range("B296").Value = 57802
$ws.Range("D296").Value = 143.48
$ws.Range("E296").Value = 162.71
$ws.Range("F296").Value = -79
$ws.Range("G296").Value = -11334.92
$ws.Range("B299").Value = 63510
$ws.Range("D299").Value = 47.64
$ws.Range("E299").Value = 50.66
$ws.Range("F299").Value = 148
$ws.Range("G299").Value = 7050.72
$ws.Range("B300").Value = 55356
$ws.Range("D300").Value = 47.64
$ws.Range("E300").Value = 54.04
$ws.Range("F300").Value = -158
$ws.Range("G300").Value = -7527.12
$ws.Range("B315").Value = 63560
$ws.Range("D315").Value = 126.86
$ws.Range("E315").Value = 134.87
$ws.Range("F315").Value = 1
$ws.Range("G315").Value = 126.86
$ws.Range("B316").Value = 60325
$ws.Range("D316").Value = 126.86
$ws.Range("E316").Value = 151.57
$ws.Range("F316").Value = -102
$ws.Range("G316").Value = -12939.72
$ws.Range("B356").Value = 31930
$ws.Range("D356").Value = 22.42
$ws.Range("E356").Value = 26.8
$ws.Range("F356").Value = -62
$ws.Range("G356").Value = -1390.04
$ws.Range("B357").Value = 63681
$ws.Range("D357").Value = 22.42
$ws.Range("E357").Value = 23.84
$ws.Range("F357").Value = 0
$ws.Range("G357").Value = 0
$ws.Range("B420").Value = 47097
$ws.Range("D420").Value = 112.28
$ws.Range("E420").Value = 134.16
$ws.Range("F420").Value = 15
$ws.Range("G420").Value = 1684.2
$ws.Range("B421").Value = 58047
$ws.Range("D421").Value = 105.54
$ws.Range("E421").Value = 126.1
$ws.Range("F421").Value = 43
$ws.Range("G421").Value = 4538.22
$ws.Range("B465").Value = 65069
$ws.Range("D465").Value = 13.45
$ws.Range("E465").Value = 14.3
$ws.Range("F465").Value = 2
$ws.Range("G465").Value = 26.9
$ws.Range("B466").Value = 53757
$ws.Range("D466").Value = 13.45
$ws.Range("E466").Value = 16.08
$ws.Range("F466").Value = -159
$ws.Range("G466").Value = -2138.55
$ws.Range("B472").Value = 45695
$ws.Range("D472").Value = 19.73
$ws.Range("E472").Value = 23.58
$ws.Range("F472").Value = -36
$ws.Range("G472").Value = -710.28
$ws.Range("B473").Value = 64915
$ws.Range("D473").Value = 19.73
$ws.Range("E473").Value = 20.98
$ws.Range("F473").Value = 0
$ws.Range("G473").Value = 0
$ws.Range("B476").Value = 45706
$ws.Range("D476").Value = 19.73
$ws.Range("E476").Value = 23.58
$ws.Range("F476").Value = -202
$ws.Range("G476").Value = -3985.46
$ws.Range("B477").Value = 64922
$ws.Range("D477").Value = 19.73
$ws.Range("E477").Value = 20.98
$ws.Range("F477").Value = 136
$ws.Range("G477").Value = 2683.28
$ws.Range("B479").Value = 64927
$ws.Range("D479").Value = 16.22
$ws.Range("E479").Value = 17.26
$ws.Range("F479").Value = 217
$ws.Range("G479").Value = 3519.74
$ws.Range("B480").Value = 45718
$ws.Range("D480").Value = 16.22
$ws.Range("E480").Value = 19.38
$ws.Range("F480").Value = -294
$ws.Range("G480").Value = -4768.68
$ws.Range("B487").Value = 45702
$ws.Range("D487").Value = 26.3
$ws.Range("E487").Value = 31.43
$ws.Range("F487").Value = -215
$ws.Range("G487").Value = -5654.5
$ws.Range("B488").Value = 64919
$ws.Range("D488").Value = 26.3
$ws.Range("E488").Value = 27.97
$ws.Range("F488").Value = 124
$ws.Range("G488").Value = 3261.2
$ws.Range("B585").Value = 60025
$ws.Range("D585").Value = 32.83
$ws.Range("E585").Value = 37.22
$ws.Range("F585").Value = -98
$ws.Range("G585").Value = -3217.34
$ws.Range("B586").Value = 64833
$ws.Range("D586").Value = 32.83
$ws.Range("E586").Value = 34.9
$ws.Range("F586").Value = 96
$ws.Range("G586").Value = 3151.68
$ws.Range("B591").Value = 64836
$ws.Range("D591").Value = 98.5
$ws.Range("E591").Value = 104.71
$ws.Range("F591").Value = 3
$ws.Range("G591").Value = 295.5
$ws.Range("B592").Value = 60031
$ws.Range("D592").Value = 98.5
$ws.Range("E592").Value = 111.69
$ws.Range("F592").Value = -5
$ws.Range("G592").Value = -492.5
$ws.Range("B732").Value = 65079
$ws.Range("D732").Value = 40.87
$ws.Range("E732").Value = 43.44
$ws.Range("F732").Value = 21
$ws.Range("G732").Value = 858.27
$ws.Range("B733").Value = 65362
$ws.Range("D733").Value = 40.87
$ws.Range("E733").Value = 43.44
$ws.Range("F733").Value = 69
$ws.Range("G733").Value = 2820.03
